$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up by one.
$ws.Rows.Item(26).Delete()

# After the above delete, the row that was "SC 92" (originally row 28) is now row 27.
# Delete it too - remaining rows shift up by one more.
$ws.Rows.Item(27).Delete()

# Now rows 26-33 hold (in order): SC5, SC101, SC105, SC119, SC120, SC132, SC193, SC232
# Update column B to reflect the new missing-data pattern.
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("B29").Value = -19.5
$ws.Range("B30").Value = -19.7
$ws.Range("B31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("B33").Value = -19.5
